# Auto-generated script to apply scheduled-runner value updates to Sheets/Excalibur_Profits.xlsx
# Updates cached market-price-derived values in columns H-N across all profession sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1757.0426
$ws.Range("J17").Value = 1907.5946
$ws.Range("L17").Value = 5722.783799999999
$ws.Range("N17").Value = -6058.783799999999
$ws.Range("H113").Value = 1437.3572
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H121").Value = 782.4545000000001
$ws.Range("J121").Value = 782.4545000000001
$ws.Range("L121").Value = 2347.3635
$ws.Range("N121").Value = -5841.3635
$ws.Range("H129").Value = 1553.4
$ws.Range("I129").Value = 1215.4615
$ws.Range("K129").Value = 3646.3845
$ws.Range("M129").Value = 1353.6155
$ws.Range("H132").Value = 29378.94
$ws.Range("I132").Value = 33045.29
$ws.Range("K132").Value = 99135.87
$ws.Range("M132").Value = -96605.87
$ws.Range("H137").Value = 1117253.1
$ws.Range("I137").Value = 1543.5714
$ws.Range("J137").Value = 2036072.9
$ws.Range("K137").Value = 4630.7142
$ws.Range("L137").Value = 6108218.699999999
$ws.Range("M137").Value = -2080.7142
$ws.Range("N137").Value = -6113318.699999999
$ws.Range("H138").Value = 3077.2188
$ws.Range("J138").Value = 3483.1667
$ws.Range("L138").Value = 10449.5001
$ws.Range("N138").Value = -20729.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5666.03
$ws.Range("I32").Value = 2355.4268
$ws.Range("K32").Value = 2355.4268
$ws.Range("M32").Value = -2068.4268
$ws.Range("H61").Value = 5758.5674
$ws.Range("I61").Value = 2425.8462
$ws.Range("K61").Value = 2425.8462
$ws.Range("M61").Value = -2213.8462
$ws.Range("H74").Value = 2671.4243
$ws.Range("I74").Value = 2547.359
$ws.Range("K74").Value = 2547.359
$ws.Range("M74").Value = -1673.359
$ws.Range("H77").Value = 2671.4243
$ws.Range("I77").Value = 2547.359
$ws.Range("K77").Value = 12736.795
$ws.Range("M77").Value = -8368.795
$ws.Range("H102").Value = 2065.5862
$ws.Range("I102").Value = 1552.2963
$ws.Range("K102").Value = 1552.2963
$ws.Range("M102").Value = 69.70370000000003
$ws.Range("H112").Value = 38887
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 38887
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 38887
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -41841
$ws.Range("H136").Value = 5758.5674
$ws.Range("I136").Value = 2425.8462
$ws.Range("K136").Value = 7277.5386
$ws.Range("M136").Value = -4727.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1888.6666
$ws.Range("I20").Value = 1859.6111
$ws.Range("J20").Value = 2063
$ws.Range("K20").Value = 1859.6111
$ws.Range("L20").Value = 2063
$ws.Range("M20").Value = -1612.6111
$ws.Range("N20").Value = -2557
$ws.Range("H80").Value = 4848.625
$ws.Range("I80").Value = 7456.4287
$ws.Range("J80").Value = 2820.3333
$ws.Range("K80").Value = 7456.4287
$ws.Range("L80").Value = 2820.3333
$ws.Range("M80").Value = -6458.4287
$ws.Range("N80").Value = -4816.3333
$ws.Range("H83").Value = 4848.625
$ws.Range("I83").Value = 7456.4287
$ws.Range("J83").Value = 2820.3333
$ws.Range("K83").Value = 37282.14350000001
$ws.Range("L83").Value = 14101.6665
$ws.Range("M83").Value = -32290.14350000001
$ws.Range("N83").Value = -24085.6665
$ws.Range("H130").Value = 177000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 177000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 177000
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = -187040
$ws.Range("H134").Value = 2745.3704
$ws.Range("I134").Value = 1893.75
$ws.Range("J134").Value = 5951.4707
$ws.Range("K134").Value = 5681.25
$ws.Range("L134").Value = 17854.4121
$ws.Range("M134").Value = -3146.25
$ws.Range("N134").Value = -22924.4121

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6506.8184
$ws.Range("I99").Value = 5844
$ws.Range("J99").Value = 7302.2
$ws.Range("K99").Value = 5844
$ws.Range("L99").Value = 7302.2
$ws.Range("M99").Value = -4346
$ws.Range("N99").Value = -10298.2
$ws.Range("H126").Value = 6506.8184
$ws.Range("I126").Value = 5844
$ws.Range("J126").Value = 7302.2
$ws.Range("K126").Value = 17532
$ws.Range("L126").Value = 21906.6
$ws.Range("M126").Value = -15062
$ws.Range("N126").Value = -26846.6
$ws.Range("H132").Value = 1417.0857
$ws.Range("I132").Value = 1494.0358
$ws.Range("J132").Value = 1109.2858
$ws.Range("K132").Value = 4482.107400000001
$ws.Range("L132").Value = 3327.8574
$ws.Range("M132").Value = -1952.107400000001
$ws.Range("N132").Value = -8387.857400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 21.333334
$ws.Range("I8").Value = 21.333334
$ws.Range("K8").Value = 64.00000199999999
$ws.Range("M8").Value = 74.99999800000001
$ws.Range("H12").Value = 16
$ws.Range("I12").Value = 8.666667
$ws.Range("K12").Value = 26.000001
$ws.Range("M12").Value = 146.999999
$ws.Range("H14").Value = 154735.84
$ws.Range("I14").Value = 154735.84
$ws.Range("K14").Value = 464207.52
$ws.Range("M14").Value = -464034.52
$ws.Range("H107").Value = 1160.5
$ws.Range("J107").Value = 1421.2858
$ws.Range("L107").Value = 4263.857400000001
$ws.Range("N107").Value = -8103.857400000001
$ws.Range("H122").Value = 883.8333
$ws.Range("J122").Value = 873.6
$ws.Range("L122").Value = 7862.400000000001
$ws.Range("N122").Value = -12762.4
$ws.Range("H132").Value = 2929.5789
$ws.Range("J132").Value = 3381.6667
$ws.Range("L132").Value = 30435.0003
$ws.Range("N132").Value = -35495.0003
$ws.Range("H140").Value = 2705.842
$ws.Range("I140").Value = 2576.5557
$ws.Range("K140").Value = 7729.6671
$ws.Range("M140").Value = -2549.6671

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5965.8335
$ws.Range("I70").Value = 5900
$ws.Range("J70").Value = 5979
$ws.Range("K70").Value = 5900
$ws.Range("L70").Value = 5979
$ws.Range("M70").Value = -5630
$ws.Range("N70").Value = -6519
$ws.Range("H73").Value = 5965.8335
$ws.Range("I73").Value = 5900
$ws.Range("J73").Value = 5979
$ws.Range("K73").Value = 5900
$ws.Range("L73").Value = 5979
$ws.Range("M73").Value = -4964
$ws.Range("N73").Value = -7851
$ws.Range("H97").Value = 6895.7144
$ws.Range("I97").Value = 1000
$ws.Range("K97").Value = 1000
$ws.Range("M97").Value = -504
$ws.Range("H102").Value = 2367.9355
$ws.Range("I102").Value = 1350.5
$ws.Range("J102").Value = 5856.2856
$ws.Range("K102").Value = 1350.5
$ws.Range("L102").Value = 5856.2856
$ws.Range("M102").Value = 271.5
$ws.Range("N102").Value = -9100.285599999999
$ws.Range("H126").Value = 3586.0334
$ws.Range("I126").Value = 1873.421
$ws.Range("J126").Value = 6544.1816
$ws.Range("K126").Value = 5620.263
$ws.Range("L126").Value = 19632.5448
$ws.Range("M126").Value = -3150.263
$ws.Range("N126").Value = -24572.5448

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3536.077
$ws.Range("I40").Value = 3269.9092
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 3269.9092
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -3133.9092
$ws.Range("N40").Value = -5272
$ws.Range("H55").Value = 631.2174
$ws.Range("I55").Value = 150
$ws.Range("K55").Value = 150
$ws.Range("M55").Value = 23
$ws.Range("H61").Value = 1707.8572
$ws.Range("I61").Value = 1836.1428
$ws.Range("J61").Value = 1579.5714
$ws.Range("K61").Value = 1836.1428
$ws.Range("L61").Value = 1579.5714
$ws.Range("M61").Value = -1634.1428
$ws.Range("N61").Value = -1983.5714
$ws.Range("H113").Value = 1707.8572
$ws.Range("I113").Value = 1836.1428
$ws.Range("J113").Value = 1579.5714
$ws.Range("K113").Value = 1836.1428
$ws.Range("L113").Value = 1579.5714
$ws.Range("M113").Value = 333.8571999999999
$ws.Range("N113").Value = -5919.5714
$ws.Range("H122").Value = 41286.332
$ws.Range("I122").Value = 3704
$ws.Range("K122").Value = 11112
$ws.Range("M122").Value = -8662

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H37").Value = 20000
$ws.Range("I37").Value = 20000
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 20000
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -19797
$ws.Range("N37").ClearContents()
$ws.Range("H55").Value = 15729.5
$ws.Range("I55").Value = 29950
$ws.Range("K55").Value = 29950
$ws.Range("M55").Value = -29673
$ws.Range("H98").Value = 71893.164
$ws.Range("J98").Value = 80271.8
$ws.Range("L98").Value = 80271.8
$ws.Range("N98").Value = -86261.8

